$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F84").Value = "DataManager.ensure_symbol_coverage implemented in backend/app/data_manager.py; uses price_bars, base_timeframe, and settings.base_horizon_days for cache-aware fetch decisions."
$ws.Range("G84").Value = "implemented"

$ws.Range("F85").Value = "BacktestService.run_single_backtest and run_group_backtest now call DataManager.ensure_symbol_coverage before loading price data; backtests no longer talk to Kite/yfinance directly."
$ws.Range("G85").Value = "implemented"

$ws.Range("F86").Value = "Regression coverage added via backend/tests/test_backtests_api.py and backend/tests/test_data_fetch_api.py so that backtests run without prior manual fetches, using synthetic sources for tests."
$ws.Range("G86").Value = "implemented"

$ws.Range("F87").Value = 'Data page now has a single "Save for backtesting (cache mode)" checkbox; in cache mode the fetch payload is adjusted to use cache-friendly timeframe/duration defaults.'
$ws.Range("G87").Value = "implemented"

$ws.Range("F88").Value = "Coverage Summary table extended with Days and BT-ready (3Y) columns, computed from created_at/start/end to indicate rows that fully cover the base horizon."
$ws.Range("G88").Value = "implemented"
